$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.152.04"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.911.07"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "362.06"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.52"
$ws.Range("E6").Value = "  -4.57%  "
$ws.Range("E7").Value = "  -4.69%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -6.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.85"
$ws.Range("E10").Value = "  -4.73%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("E12").Value = "  -3.90%  "
$ws.Range("E13").Value = "  -5.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.362.79"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.32"
$ws.Range("E15").Value = "  -5.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.904.96"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.952"
$ws.Range("E17").Value = "  -2.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.038.35"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.29"
$ws.Range("E19").Value = "  -3.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.21"
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.97"
$ws.Range("E21").Value = "  -6.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0945"
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.15"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "258.43"
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.67"
$ws.Range("E25").Value = "  -4.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.171"
$ws.Range("E26").Value = "  -6.06%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.83"
$ws.Range("E28").Value = "  -3.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.11"
$ws.Range("E29").Value = "  -6.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.106"
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.90"
$ws.Range("E32").Value = "  -5.23%  "
$ws.Range("E33").Value = "  -2.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.82"
$ws.Range("E34").Value = "  -6.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.51"
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0419"
$ws.Range("E37").Value = "  -5.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.80"
$ws.Range("E38").Value = "  +3.11%  "
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.88"
$ws.Range("E40").Value = "  -6.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.86"
$ws.Range("E41").Value = "  -6.44%  "
$ws.Range("E42").Value = "  -4.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.32"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.14"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.12"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.059.67"
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.18"
$ws.Range("E47").Value = "  -7.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.26"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.188.78"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.235"
$ws.Range("E50").Value = "  -6.36%  "
$ws.Range("E51").Value = "  -8.10%  "
